$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the refreshed cryptos data.
# For column D we force a Text number format before assigning the value so that
# numeric-looking strings (e.g. "1.00", "0.0000289") are preserved exactly as text
# instead of being auto-converted by Excel into floating point numbers. The style
# is then reset back to "Normal" so no new cell style / number format is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.416.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.616.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.609.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.665"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000289"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("E14").Value = "  +5.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.195.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.617.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.389.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.44%  "

$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("E21").Value = "  +4.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.82%  "

$ws.Range("E25").Value = "  +2.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.01"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.07%  "

$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.89%  "

$ws.Range("E30").Value = "  +5.53%  "

$ws.Range("E31").Value = "  +8.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "639.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.120"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0823"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.02%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.307.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.58%  "

$ws.Range("E45").Value = "  +5.41%  "

$ws.Range("E46").Value = "  +2.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "

$ws.Range("E48").Value = "  +5.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
